# Fruta / hortaliza, semanal
# Insert a new price-report row above row 85, shifting the existing
# rows 85-154 down to 86-155, and populate the newly inserted row 85
# with the new daily entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 85 (pushes rows 85..154 to 86..155).
$ws.Rows("85:85").Insert()

# Populate the new row 85 with the newly-reported values.
$ws.Range("A85").Value = 8
$ws.Range("B85").Value = "Terminal La Palmera de La Serena"
$ws.Range("C85").Value = "Coquimbo"
$ws.Range("D85").Value = 45068
$ws.Range("E85").Value = 4
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100109
$ws.Range("H85").Value = "Uva"
$ws.Range("I85").Value = 100109001
$ws.Range("J85").Value = "Uva"
$ws.Range("K85").Value = "Autumn Royal"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 240
$ws.Range("N85").Value = 8000
$ws.Range("O85").Value = 9000
$ws.Range("P85").Value = 8500
$ws.Range("Q85").Value = "$/bandeja 18 kilos"
$ws.Range("R85").Value = "Provincia de Limarí"
$ws.Range("S85").Value = 472
$ws.Range("T85").Value = 18
